# Update paypal chrome download dialog missing extension
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet2: fill in D2:D4 with the "False" status value (matches D1 header "Status").
# Column C already holds that same literal text "False" (t="s") in rows 2-4, so
# copy/paste-special (values) from C into D to replicate the text cell faithfully
# instead of letting a direct Value assignment coerce "False" into a Boolean.
$ws2.Range("C2:C4").Copy() | Out-Null
$ws2.Range("D2:D4").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Update the active selections on each sheet to match the saved view state
$ws1.Range("A4:XFD14").Select() | Out-Null
$ws2.Range("D2:D4").Select() | Out-Null

# Leave Sheet1 as the active / tab-selected sheet
$ws1.Activate()
$ws1.Range("A4:XFD14").Select() | Out-Null
